$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update NS row (row 7): NB_BLESSES value changes from 0 to 4
$ws.Range("B7").Value = 4

# Add new row for R (Rugby)
$ws.Range("A8").Value = "R"
$ws.Range("B8").Value = 7
$ws.Range("C8").Value = 14

# Rename header B1 from NB_BLESSURES to NB_BLESSES
$ws.Range("B1").Value = "NB_BLESSES"

# Add new row for BS (Baseball)
$ws.Range("A9").Value = "BS"
$ws.Range("B9").Value = 17
$ws.Range("C9").Value = 21

# Fill the BLESS_POURC formula down through the new rows
$ws.Range("D3:D9").Formula = "=(B3/C3)*100"

# Update selection to D9 to match final state
$ws.Range("D9").Select()

$wb.Save()
